# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values in column F ("dSF") for specific rows, per the repulled data.
$updates = @{
    2  = 0
    3  = -2
    13 = 0
    15 = -1
    18 = 5
    27 = -9
    28 = -2
    29 = -4
    31 = 0
    44 = 2
    48 = 1
    53 = -1
    55 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
